$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.752.55"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.887.51"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'0.7929"
$ws.Range("E5").Value = "  -5.50%  "
$ws.Range("D6").Value = "'241.69"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.3168"
$ws.Range("D9").Value = "'25.41"
$ws.Range("E9").Value = "  -5.40%  "
$ws.Range("D10").Value = "'0.07001"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "'0.08045"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "'0.7660"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").Value = "1.889.28"
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "'5.302"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "'92.20"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "29.755.29"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "'13.86"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").Value = "'5.943"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "'243.41"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "'0.000007687"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'8.131"
$ws.Range("E22").Value = "  +15.37%  "
$ws.Range("D23").Value = "2.136.76"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'0.1676"
$ws.Range("E25").Value = "  +4.38%  "
$ws.Range("D26").Value = "'9.303"
$ws.Range("D27").Value = "'164.76"
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("D28").Value = "'18.64"
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("D29").Value = "'2.053"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").Value = "'1.533"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").Value = "'4.384"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").Value = "'0.05665"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").Value = "'4.044"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").Value = "'1.262"
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("D36").Value = "'0.7341"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "'0.9989"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "'2.638"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").Value = "'0.01910"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").Value = "'2.765"
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("D41").Value = "'0.4404"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("D42").Value = "'72.43"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "'5.825"
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").Value = "'1.0000"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "'0.8366"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "'102.71"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "1.021.54"
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("D48").Value = "'1.863"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").Value = "'9.872"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "'7.426"
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("D51").Value = "2.018.17"
$ws.Range("E51").Value = "  -2.12%  "
